$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update harvester (column B) for all data rows, and set the new
# experimentDesign value (column D) that was previously blank.
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"
}

# Reflect the selection left behind by the author's edit (D2:D16, active cell D2)
$ws.Range("D2:D16").Select()
